$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "folder_type" style tag columns (B = folder_id-like tag,
# C = "ofs" source marker, D = ofs sub-category) for each dataset row,
# based on the dataset type classification requested in the commit.

$ws.Range("B2").Value  = "17-politique"

$ws.Range("B3").Value  = "11-mobilite"

$ws.Range("B4").Value  = "14-sante"

$ws.Range("B5").Value  = "02-espace"

$ws.Range("B6").Value  = "01-pop"
$ws.Range("C6").Value  = "ofs"
$ws.Range("D6").Value  = "ofs-div-pop"

$ws.Range("B7").Value  = "04-economie"

$ws.Range("B8").Value  = "14-sante"

$ws.Range("B9").Value  = "10-tourisme"
$ws.Range("C9").Value  = "ofs"
$ws.Range("D9").Value  = "ofs-tourisme"

$ws.Range("B10").Value = "02-espace"

$ws.Range("B11").Value = "11-mobilite"

$ws.Range("B12").Value = "01-pop"

$ws.Range("B13").Value = "03-travail"

$ws.Range("B14").Value = "14-sante"

$ws.Range("B15").Value = "11-mobilite"

$ws.Range("B16").Value = "06-industrie"

$ws.Range("B17").Value = "12-monnaie"

$ws.Range("B18").Value = "08-energie"

$ws.Range("B19").Value = "04-economie"
$ws.Range("C19").Value = "ofs"
$ws.Range("D19").Value = "ofs-salaire"

$ws.Range("B20").Value = "03-travail"
$ws.Range("C20").Value = "ofs"
$ws.Range("D20").Value = "ofs-travail"

# Widen column E to fit the new content layout.
$ws.Columns.Item(5).ColumnWidth = 27.33

# Update the selected/active cell in the frozen-pane view.
$ws.Range("D8").Select()
